# premier league - atualização de dados
# atualizando rodada 5 do sabado - 22 de setembro
#
# Updates the "cartoes" (cards) analysis sheet: several clubs played an
# extra match (rodada 5) on Saturday, 22 de setembro. This bumps their
# total matches played (B) and either their home matches (C) or away
# matches (F) by one, along with the related yellow-card totals
# (D/G/I) and the derived averages (E/H/J).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Map of cell -> new value, one block per updated club/row.
$updates = @{
    # Aston Villa (row 3) - played the extra match at home
    "B3" = 5
    "C3" = 3
    "D3" = 5
    "E3" = 1.666666666666667
    "I3" = 12
    "J3" = 2.4

    # Bournemouth (row 4) - played the extra match away
    "B4" = 5
    "F4" = 3
    "G4" = 8
    "H4" = 2.666666666666667
    "I4" = 16
    "J4" = 3.2

    # Brentford (row 5) - played the extra match away
    "B5" = 5
    "F5" = 3
    "G5" = 5
    "H5" = 1.666666666666667
    "I5" = 8
    "J5" = 1.6

    # Chelsea (row 7) - played the extra match away
    "B7" = 5
    "F7" = 3
    "G7" = 13
    "H7" = 4.333333333333333
    "I7" = 18
    "J7" = 3.6

    # Crystal Palace (row 9) - played the extra match at home
    "B9" = 5
    "C9" = 3
    "D9" = 4
    "E9" = 1.333333333333333
    "I9" = 11
    "J9" = 2.2

    # Everton (row 10) - played the extra match away
    "B10" = 5
    "F10" = 3
    "G10" = 6
    "I10" = 9
    "J10" = 1.8

    # Fulham (row 11) - played the extra match at home
    "B11" = 5
    "C11" = 3
    "D11" = 10
    "E11" = 3.333333333333333
    "I11" = 16
    "J11" = 3.2

    # Ipswich Town (row 12) - played the extra match away
    "B12" = 5
    "F12" = 3
    "G12" = 10
    "H12" = 3.333333333333333
    "I12" = 15
    "J12" = 3

    # Leicester City (row 13) - played the extra match at home
    "B13" = 5
    "C13" = 3
    "D13" = 6
    "E13" = 2
    "I13" = 11
    "J13" = 2.2

    # Liverpool (row 14) - played the extra match at home
    "B14" = 5
    "C14" = 3
    "D14" = 7
    "E14" = 2.333333333333333
    "I14" = 9
    "J14" = 1.8

    # Manchester United (row 15) - played the extra match away
    "B15" = 5
    "F15" = 3
    "G15" = 8
    "H15" = 2.666666666666667
    "I15" = 14
    "J15" = 2.8

    # Newcastle (row 16) - played the extra match away
    "B16" = 5
    "F16" = 3
    "H16" = 1.666666666666667
    "J16" = 2.2

    # Southampton (row 18) - played the extra match at home
    "B18" = 5
    "C18" = 3
    "D18" = 8
    "E18" = 2.666666666666667
    "I18" = 13
    "J18" = 2.6

    # Tottenham (row 19) - played the extra match at home
    "B19" = 5
    "C19" = 3
    "D19" = 8
    "E19" = 2.666666666666667
    "I19" = 13
    "J19" = 2.6

    # West Ham (row 20) - played the extra match at home
    "B20" = 5
    "C20" = 3
    "D20" = 9
    "E20" = 3
    "I20" = 13
    "J20" = 2.6

    # Wolves (row 21) - played the extra match away
    "B21" = 5
    "F21" = 3
    "G21" = 12
    "H21" = 4
    "I21" = 17
    "J21" = 3.4
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
